# Generate Report for Handoff
# Update the localization status for the ed4a6625 file (row 3) to reflect that
# a new handoff package ("Ready for handoff") has just been generated.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-09-05 10:16:16"
$overview.Columns.Item(5).ColumnWidth = 16.25
$overview.Columns.Item(6).ColumnWidth = 16.25

# --- zh-cn sheet ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("E3").Value = "mt"
$zhcn.Range("H3").Value = "2016-09-05 10:16:11"
$zhcn.Columns.Item(3).ColumnWidth = 16.25

# --- de-de sheet ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("E3").Value = "mt"
$dede.Range("H3").Value = "2016-09-05 10:16:16"
$dede.Columns.Item(3).ColumnWidth = 16.25
